# Tournament deck breakdown - update hash csv - update meta archetype - new module stat_helper
#
# Applies the archetype-table refresh to Sheet1 of FortuneHandMeta.xlsx:
#  - Row 4  (was "Roach-Terror Forest"): renamed/re-identified as "Roach Accel Forest"
#  - Row 5  (was "Pure Roach Forest"):   renamed to "Roach Natura Forest"
#  - Row 7  (Aggro Rally Sword):         anti-identifier swapped to a new card, 3rd slot cleared
#  - Row 15 (Jerva Dragon):              renamed "Jerva/Evo Dragon", identifier 2 updated
#  - Row 21 (Control Blood):             identifier 2 + anti-identifier updated
#  - Row 25 (Amulet Haven):              renamed "Sofina/Amulet Haven", identifier updated
#  - Cosmetic: column C width, and the active selection/scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 : Roach-Terror Forest -> Roach Accel Forest ---------------------
$ws.Range("A4").Value = "Roach Accel Forest"
$ws.Range("C4").Value = "Resolve of the Fallen"
$ws.Range("D4").Value = "6-SmQ"
$ws.Range("G4").Value = "Fertile Aether"
$ws.Range("H4").Value = "6pMJg"

# --- Row 5 : Pure Roach Forest -> Roach Natura Forest -----------------------
$ws.Range("A5").Value = "Roach Natura Forest"

# --- Row 7 : Aggro Rally Sword anti-identifier swap -------------------------
$ws.Range("C7").Value = "Ernesta, Weapons Hawker"
$ws.Range("G7").Value = "None"
$ws.Range("H7").Value = "None"

# --- Row 15 : Jerva Dragon -> Jerva/Evo Dragon ------------------------------
$ws.Range("A15").Value = "Jerva/Evo Dragon"
$ws.Range("E15").Value = "Jerva, Wyrm Transcendent"
$ws.Range("F15").Value = "6yB_6"

# --- Row 21 : Control Blood identifier 2 / anti-identifier update ----------
$ws.Range("E21").Value = "Nerea, Beast Empress"
$ws.Range("F21").Value = "6yypo"
$ws.Range("G21").Value = "Permafrost Behemoth"
$ws.Range("H21").Value = "6v8h6"

# --- Row 25 : Amulet Haven -> Sofina/Amulet Haven ---------------------------
$ws.Range("A25").Value = "Sofina/Amulet Haven"
$ws.Range("C25").Value = "VIII. Sofina, Strength"
$ws.Range("D25").Value = "719NI"

# --- Cosmetic formatting / view state ---------------------------------------
$ws.Columns.Item(3).ColumnWidth = 26.85546875

$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("G9").Select()
